$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 2714
$ws1.Range("F13").Value = 55
$ws1.Range("F16").Value = 361
$ws1.Range("F18").Value = 2137
$ws1.Range("F20").Value = 720

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4218

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 219

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 2714
$ws4.Range("F18").Value = 361
$ws4.Range("F22").Value = 2137
$ws4.Range("F24").Value = 720
$ws4.Range("F34").Value = 219
